$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the cell to stay a text string even when the text looks numeric
    # (mirrors the inline-string / non-coerced representation in the source file).
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.923.62"
$ws.Range("E2").Value = "  +0.90%  "

Set-TextValue $ws.Range("D3") "1.634.49"
$ws.Range("E3").Value = "  +1.59%  "

$ws.Range("E4").Value = "  +0.46%  "

Set-TextValue $ws.Range("D5") "214.77"
$ws.Range("E5").Value = "  +1.15%  "

Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.46%  "

Set-TextValue $ws.Range("D8") "28.79"
$ws.Range("E8").Value = "  -0.43%  "

$ws.Range("E9").Value = "  -0.21%  "

$ws.Range("E10").Value = "  +0.37%  "

$ws.Range("E11").Value = "  -0.33%  "

Set-TextValue $ws.Range("D12") "1.870.78"
$ws.Range("E12").Value = "  +1.74%  "

Set-TextValue $ws.Range("D13") "1.620.13"
$ws.Range("E13").Value = "  +0.74%  "

Set-TextValue $ws.Range("D14") "0.563"
$ws.Range("E14").Value = "  -0.11%  "

Set-TextValue $ws.Range("D15") "9.30"
$ws.Range("E15").Value = "  +8.77%  "

Set-TextValue $ws.Range("D16") "29.954.07"
$ws.Range("E16").Value = "  +0.94%  "

Set-TextValue $ws.Range("D17") "3.85"
$ws.Range("E17").Value = "  +0.32%  "

Set-TextValue $ws.Range("D18") "64.22"
$ws.Range("E18").Value = "  -0.62%  "

Set-TextValue $ws.Range("D19") "241.82"
$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("E20").Value = "  -0.10%  "

$ws.Range("E21").Value = "  +0.31%  "

$ws.Range("E22").Value = "  +1.87%  "

Set-TextValue $ws.Range("D23") "9.79"
$ws.Range("E23").Value = "  +2.27%  "

Set-TextValue $ws.Range("D24") "2.16"
$ws.Range("E24").Value = "  +2.85%  "

Set-TextValue $ws.Range("D25") "157.98"
$ws.Range("E25").Value = "  +0.87%  "

Set-TextValue $ws.Range("D26") "15.47"
$ws.Range("E26").Value = "  -0.56%  "

$ws.Range("E27").Value = "  +0.28%  "

Set-TextValue $ws.Range("D28") "6.58"
$ws.Range("E28").Value = "  +0.62%  "

$ws.Range("E29").Value = "  +0.43%  "

Set-TextValue $ws.Range("D30") "0.0490"
$ws.Range("E30").Value = "  +1.86%  "

$ws.Range("E33").Value = "  +0.08%  "

Set-TextValue $ws.Range("D34") "1.428.03"
$ws.Range("E34").Value = "  +0.20%  "

$ws.Range("E35").Value = "  +4.94%  "

$ws.Range("E36").Value = "  -1.11%  "

Set-TextValue $ws.Range("D37") "2.75"
$ws.Range("E37").Value = "  -4.50%  "

$ws.Range("E38").Value = "  +0.43%  "

$ws.Range("E39").Value = "  +0.40%  "

Set-TextValue $ws.Range("D40") "75.62"
$ws.Range("E40").Value = "  +10.68%  "

Set-TextValue $ws.Range("D41") "0.555"
$ws.Range("E41").Value = "  +0.36%  "

Set-TextValue $ws.Range("D42") "0.0503"
$ws.Range("E42").Value = "  +1.09%  "

$ws.Range("E43").Value = "  +0.83%  "

Set-TextValue $ws.Range("D44") "0.828"
$ws.Range("E44").Value = "  +0.49%  "

$ws.Range("E45").Value = "  +0.44%  "

$ws.Range("E46").Value = "  +0.64%  "

Set-TextValue $ws.Range("D47") "51.00"
$ws.Range("E47").Value = "  -6.03%  "

Set-TextValue $ws.Range("D48") "1.777.07"
$ws.Range("E48").Value = "  +1.69%  "

$ws.Range("E49").Value = "  -1.41%  "

Set-TextValue $ws.Range("D50") "90.51"
$ws.Range("E50").Value = "  +3.76%  "

Set-TextValue $ws.Range("D51") "0.0₆0112"
$ws.Range("E51").Value = "  +9.41%  "

# Rows 31/32: Filecoin <-> PancakeSwap data updated (not a literal cell swap -- new values differ slightly)
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D31") "1.10"
$ws.Range("E31").Value = "  +3.37%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D32") "3.38"
$ws.Range("E32").Value = "  +3.23%  "
